# "Fix to resolve end point" -- scrub the sample/demo avail rows so the
# template no longer ships with leftover dummy test data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the placeholder titles to generic sample titles.
$ws.Range("A3").Value = "Sample Title 1"
$ws.Range("A4").Value = "Sample, Title 2"

# Clear the sample Episode Price / Episode SR Price / Pre-Order Date /
# Suppression End Date values for both example rows.
$ws.Range("Q3:T3").ClearContents()
$ws.Range("Q4:T4").ClearContents()

# Clear the sample Local Data Rating / Local Data # of Episodes /
# Local Edit Required values for both example rows.
$ws.Range("W3:Y3").ClearContents()
$ws.Range("W4:Y4").ClearContents()

# Update the sample Vendor ID on row 4.
$ws.Range("AB4").Value = 123456

# Leave the selection on the second example row, matching the saved view.
$ws.Range("A4").Select()
